# Apply the two content edits described by the diff:
#  1. Slide 5 (content placeholder, "Fin." paragraph that follows): the end date
#     "Miércoles 25 de Noviembre del 2015." becomes "Jueves 30 de Julio del " +
#     "2015." (split across two runs).
#  2. Slide 9 (content placeholder, last bullet): "65 días en total del proyecto"
#     becomes "42.63 días en total del proyecto", re-split into three runs:
#     "42.63 " / "días " / "en total del proyecto".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 5: end-date paragraph
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(1)
$endPara = $sh5.TextFrame.TextRange.Paragraphs(5)

# Replace the whole paragraph text first …
$endPara.Text = "Jueves 30 de Julio del 2015."

# … then force the run boundary right after "del " so the text ends up split
# into "Jueves 30 de Julio del " + "2015." as two separate runs.
$firstPart = $endPara.Characters(1, 23)
$firstPart.Text = "Jueves 30 de Julio del "

# ---------------------------------------------------------------------------
# Slide 9: "días en total del proyecto" paragraph
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(1)
$daysPara = $sh9.TextFrame.TextRange.Paragraphs(6)

# "65" -> "42.63"
$numberPart = $daysPara.Characters(1, 2)
$numberPart.Text = "42.63"

# Merge "42.63" + " " into a single run "42.63 "
$firstRun = $daysPara.Characters(1, 6)
$firstRun.Text = "42.63 "

# Merge "d" + "ías " into a single run "días "
$secondRun = $daysPara.Characters(7, 5)
$secondRun.Text = "días "
